# Regenerate save_data column G ("K") with newly computed values.
# (regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 2
    3  = 1
    4  = 1
    5  = 1
    6  = 2
    7  = 1
    8  = 2
    9  = 1
    10 = 1
    11 = 2
    12 = 1
    13 = 3
    14 = 0
    15 = 1
    16 = 1
    17 = 2
    18 = 1
    19 = 1
    20 = 1
    21 = 1
    22 = 0
    23 = 1
    24 = 1
    25 = 1
    26 = 1
    27 = 2
    28 = 1
    29 = 0
    31 = 4
    32 = 3
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
